$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update review text (column A) and sentiment (column B) for rows 2-53
$data = @(
    ,@(2, "i dont liek teh product,", "NEGATIVE")
    ,@(3, "i love it , its wonderful , it exceeded my expectations", "POSITIVE")
    ,@(4, "Battery life could be better.", "NEGATIVE")
    ,@(5, "Excellent quality, highly recommend!", "POSITIVE")
    ,@(6, "Durable and worth the money.", "POSITIVE")
    ,@(7, "Excellent quality, highly recommend!", "POSITIVE")
    ,@(8, "Would definitely buy again.", "POSITIVE")
    ,@(9, "Battery life could be better.", "NEGATIVE")
    ,@(10, "Would definitely buy again.", "POSITIVE")
    ,@(11, "Not as expected, a bit disappointed.", "NEGATIVE")
    ,@(12, "Amazing performance and design.", "POSITIVE")
    ,@(13, "Fast delivery and good packaging.", "POSITIVE")
    ,@(14, "Battery life could be better.", "NEGATIVE")
    ,@(15, "Good value for the price.", "POSITIVE")
    ,@(16, "Exceeded my expectations!", "POSITIVE")
    ,@(17, "Very comfortable to use.", "POSITIVE")
    ,@(18, "Durable and worth the money.", "POSITIVE")
    ,@(19, "Very comfortable to use.", "POSITIVE")
    ,@(20, "Amazing performance and design.", "POSITIVE")
    ,@(21, "Very comfortable to use.", "POSITIVE")
    ,@(22, "Would definitely buy again.", "POSITIVE")
    ,@(23, "Not as expected, a bit disappointed.", "NEGATIVE")
    ,@(24, "Durable and worth the money.", "POSITIVE")
    ,@(25, "Good value for the price.", "POSITIVE")
    ,@(26, "Fast delivery and good packaging.", "POSITIVE")
    ,@(27, "Battery life could be better.", "NEGATIVE")
    ,@(28, "Fast delivery and good packaging.", "POSITIVE")
    ,@(29, "Amazing performance and design.", "POSITIVE")
    ,@(30, "Durable and worth the money.", "POSITIVE")
    ,@(31, "Durable and worth the money.", "POSITIVE")
    ,@(32, "Excellent quality, highly recommend!", "POSITIVE")
    ,@(33, "Very comfortable to use.", "POSITIVE")
    ,@(34, "Excellent quality, highly recommend!", "POSITIVE")
    ,@(35, "Very comfortable to use.", "POSITIVE")
    ,@(36, "Excellent quality, highly recommend!", "POSITIVE")
    ,@(37, "Battery life could be better.", "NEGATIVE")
    ,@(38, "Amazing performance and design.", "POSITIVE")
    ,@(39, "Exceeded my expectations!", "POSITIVE")
    ,@(40, "Would definitely buy again.", "POSITIVE")
    ,@(41, "Durable and worth the money.", "POSITIVE")
    ,@(42, "Fast delivery and good packaging.", "POSITIVE")
    ,@(43, "Excellent quality, highly recommend!", "POSITIVE")
    ,@(44, "Good value for the price.", "POSITIVE")
    ,@(45, "Very comfortable to use.", "POSITIVE")
    ,@(46, "Excellent quality, highly recommend!", "POSITIVE")
    ,@(47, "Not as expected, a bit disappointed.", "NEGATIVE")
    ,@(48, "Durable and worth the money.", "POSITIVE")
    ,@(49, "Exceeded my expectations!", "POSITIVE")
    ,@(50, "Would definitely buy again.", "POSITIVE")
    ,@(51, "Very comfortable to use.", "POSITIVE")
    ,@(52, "Very comfortable to use.", "POSITIVE")
    ,@(53, "Durable and worth the money.", "POSITIVE")
)

foreach ($item in $data) {
    $r = $item[0]
    $ws.Cells.Item($r, 1).Value = $item[1]
    $ws.Cells.Item($r, 2).Value = $item[2]
}

# Remove rows 54-59 which are no longer part of the dataset
$ws.Range("A54:B59").EntireRow.Delete() | Out-Null
